$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 257, shifting existing rows 257:275 down to 258:276
$ws.Rows.Item(257).Insert()

# Populate the newly inserted row 257 with the new weekly record
$ws.Cells.Item(257, 1).Value = 1
$ws.Cells.Item(257, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(257, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(257, 4).Value = 44783
$ws.Cells.Item(257, 5).Value = 15
$ws.Cells.Item(257, 6).Value = "Fruta"
$ws.Cells.Item(257, 7).Value = 100108
$ws.Cells.Item(257, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(257, 9).Value = 100108006
$ws.Cells.Item(257, 10).Value = "Plátano"
$ws.Cells.Item(257, 11).Value = "Sin especificar"
$ws.Cells.Item(257, 12).Value = "Pintón"
$ws.Cells.Item(257, 13).Value = 120
$ws.Cells.Item(257, 14).Value = 21000
$ws.Cells.Item(257, 15).Value = 22000
$ws.Cells.Item(257, 16).Value = 21500
$ws.Cells.Item(257, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(257, 18).Value = "Ecuador"
$ws.Cells.Item(257, 19).Value = 1075
$ws.Cells.Item(257, 20).Value = 20
